$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 21 (day 21) - first new study log entry
$ws.Range("A5").Value = 21
$ws.Range("B5").Value = "8：30-10：15"
$ws.Range("C5").Value = "for循环 while循环 dowhile循环 死循环"
$ws.Range("D5").Value = "2：：03-3:38"
$ws.Range("E5").Value = "循环嵌套 break continue 生成随机数"

# Row 21 (day 21) - second new study log entry
$ws.Range("A6").Value = 21
$ws.Range("B6").Value = "10：17-10：46"
$ws.Range("C6").Value = "C语言：一些基本概念吧"
$ws.Range("D6").Value = "4:05-5：20"
$ws.Range("E6").Value = "C语言一些基本概念  高数 映射"

# Row 22 (day 22)
$ws.Range("A7").Value = 22
$ws.Range("B7").Value = "9：20-10：41"
$ws.Range("C7").Value = "函数"

# Row 24 (day 24)
$ws.Range("A8").Value = 24
$ws.Range("B8").Value = "9：00-10：30"
$ws.Range("C8").Value = "函数的一些性质 数列的极限（没学完）"
$ws.Range("D8").Value = "7：00-8：34"
$ws.Range("E8").Value = "翻了一下书 去翻了一些视频（算是vs的补充内容？还有一些算是一些程序员的基本修养？）"
$ws.Range("G8").Value = "emmm感觉还算有用吧，这两天就相当于构建一下环境，学习了vs，明天算是正式开始C语言吧"

# Widen column G to fit the new long note in G8
$ws.Columns.Item(7).ColumnWidth = 80

# Match the workbook's final selection/active cell
$ws.Range("G8").Select()
